$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B3 date value (2017-01-01 -> 2018-01-01)
$ws.Range("B3").Value = 43101

# Add new row 4 data
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 43108
$ws.Range("C4").Value = 16.2
$ws.Range("D4").Value = 102.5
$ws.Range("E4").Value = 225.9
$ws.Range("F4").Formula = "=E4-E3"
$ws.Range("G4").Formula = "=E4-210"

# Update selection to match diff
$ws.Range("F6").Select()
